$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9668371191043447
$ws.Range("D2").Value = 0.9837685289672864
$ws.Range("E2").Value = 0.4838891568136745
$ws.Range("F2").Value = 0.2147858542782506
$ws.Range("G2").Value = 9.712481265923016
$ws.Range("H2").Value = 18.21067843207808
$ws.Range("I2").Value = 1.031688863338307
$ws.Range("J2").Value = -0.4631657018823567
$ws.Range("C3").Value = 0.729995635389626
$ws.Range("D3").Value = 0.8546163145052088
$ws.Range("E3").Value = 1.380718541818912
$ws.Range("F3").Value = 1.06316865939333
$ws.Range("G3").Value = 48.07581822170226
$ws.Range("H3").Value = 51.9619442101981
$ws.Range("I3").Value = 1.021084093618813
$ws.Range("J3").Value = -0.3165895842406616
$ws.Range("C4").Value = 0.6975743901747746
$ws.Range("D4").Value = 0.8428913512719718
$ws.Range("E4").Value = 1.155852661010008
$ws.Range("F4").Value = 0.9032250202491747
$ws.Range("G4").Value = 52.61096343483078
$ws.Range("H4").Value = 54.99323684101759
$ws.Range("I4").Value = 0.8817721611952848
$ws.Range("J4").Value = 1.670687758040486
$ws.Range("C5").Value = 0.8640955775124071
$ws.Range("D5").Value = 0.9304685044350246
$ws.Range("E5").Value = 0.2162582362440504
$ws.Range("F5").Value = 0.01978462785931484
$ws.Range("G5").Value = 4.371814713925685
$ws.Range("H5").Value = 36.86521700568071
$ws.Range("I5").Value = 1.036083042583906
$ws.Range("J5").Value = -0.02526949527042643
$ws.Range("C6").Value = 0.4691420534716806
$ws.Range("D6").Value = 0.6876019327991751
$ws.Range("E6").Value = 0.4274103388672515
$ws.Range("F6").Value = 0.2517300785272529
$ws.Range("G6").Value = 55.62486537875293
$ws.Range("H6").Value = 72.85999907550915
$ws.Range("I6").Value = 0.9251070200291711
$ws.Range("J6").Value = 0.09680413821106781
$ws.Range("C7").Value = 0.6924979152275967
$ws.Range("D7").Value = 0.83425667698703
$ws.Range("E7").Value = 0.2422689495061176
$ws.Range("F7").Value = 0.1786910929010296
$ws.Range("G7").Value = 50.51571238869423
$ws.Range("H7").Value = 55.45287050932559
$ws.Range("I7").Value = 0.9339060197052601
$ws.Range("J7").Value = 0.07248169257056181
$ws.Range("C8").Value = 0.9944182798220433
$ws.Range("D8").Value = 0.9972885540004741
$ws.Range("E8").Value = 0.02308284357410146
$ws.Range("F8").Value = 0.003082209729107735
$ws.Range("G8").Value = 1.282890873152805
$ws.Range("H8").Value = 7.471091070223087
$ws.Range("I8").Value = 1.01119134512855
$ws.Range("J8").Value = -0.03464381070903011
$ws.Range("C9").Value = 0.459896433987896
$ws.Range("D9").Value = 0.6781705771256321
$ws.Range("E9").Value = 0.2270616557772615
$ws.Range("F9").Value = 0.1331900484210723
$ws.Range("G9").Value = 55.43694703852557
$ws.Range("H9").Value = 73.4917387202197
$ws.Range("I9").Value = 1.006071643761663
$ws.Range("J9").Value = -0.01946869070319712
$ws.Range("C10").Value = 0.7043874049974126
$ws.Range("D10").Value = 0.8404146108051278
$ws.Range("E10").Value = 0.1423039584866977
$ws.Range("F10").Value = 0.1039128710724974
$ws.Range("G10").Value = 44.33142963843744
$ws.Range("H10").Value = 54.37026715058402
$ws.Range("I10").Value = 0.963637386661083
$ws.Range("J10").Value = 0.1275200198976481
$ws.Range("C11").Value = 0.2208260315990973
$ws.Range("D11").Value = 0.7126494172498229
$ws.Range("E11").Value = 90.71738856043673
$ws.Range("F11").Value = 71.81290370180366
$ws.Range("G11").Value = 86.23904606304276
$ws.Range("H11").Value = 88.27083144509871
$ws.Range("I11").Value = 3.946403264013481
$ws.Range("J11").Value = -1486.043385673669
$ws.Range("C12").Value = 0.176246253369531
$ws.Range("D12").Value = 0.6257392525354241
$ws.Range("E12").Value = 93.27645323135607
$ws.Range("F12").Value = 74.93486195523286
$ws.Range("G12").Value = 89.9881592132692
$ws.Range("H12").Value = 90.76088070476558
$ws.Range("I12").Value = 3.768034773461992
$ws.Range("J12").Value = -1395.519948651262
$ws.Range("C13").Value = 0.2539318915221795
$ws.Range("D13").Value = 0.6869378455299529
$ws.Range("E13").Value = 66.72608544297137
$ws.Range("F13").Value = 54.29271384887015
$ws.Range("G13").Value = 85.54803170610228
$ws.Range("H13").Value = 86.37523420968654
$ws.Range("I13").Value = 3.061158717118162
$ws.Range("J13").Value = -1050.102832176017
$ws.Range("C14").Value = 0.8397624241164194
$ws.Range("D14").Value = 0.9201185133622614
$ws.Range("E14").Value = 1.287523200635957
$ws.Range("F14").Value = 0.2736037563924894
$ws.Range("G14").Value = 11.08066363821167
$ws.Range("H14").Value = 40.02968596973759
$ws.Range("I14").Value = 1.09756569867533
$ws.Range("J14").Value = -1.452002709852982
$ws.Range("C15").Value = 0.6575307647631854
$ws.Range("D15").Value = 0.8149389752724198
$ws.Range("E15").Value = 1.88227755017504
$ws.Range("F15").Value = 1.102285999917759
$ws.Range("G15").Value = 44.64142071455053
$ws.Range("H15").Value = 58.52087108347026
$ws.Range("I15").Value = 1.10861362194511
$ws.Range("J15").Value = -1.614115737407085
$ws.Range("C16").Value = 0.7879435234592955
$ws.Range("D16").Value = 0.8912800469071079
$ws.Range("E16").Value = 0.9979038192693285
$ws.Range("F16").Value = 0.776333607171133
$ws.Range("G16").Value = 43.75311151709753
$ws.Range("H16").Value = 46.04959028489878
$ws.Range("I16").Value = 0.9237086232935777
$ws.Range("J16").Value = 1.094307384111133
